$wb = $excel.ActiveWorkbook

# The same set of updates needs to be applied to both the "展览" sheet
# and the "全部类型" sheet, since they contain duplicated rows.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("G2").Value = 68
    $ws.Range("F4").Value = 81
    $ws.Range("F5").Value = 382
    $ws.Range("F6").Value = 11333
    $ws.Range("F7").Value = 664
    $ws.Range("F12").Value = 157
    $ws.Range("F19").Value = 1258
    $ws.Range("F21").Value = 891
}
